# "creating add product with external data"
#
# Changes applied to the Product_info sheet:
#   1. Categories value in B5 changes from "Electronics-others" to
#      "Electronics >> Others".
#   2. The sheet's view is scrolled back to the top (top-left visible
#      cell A1) and the active/selected cell moves from D10 to B11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Categories" text value (row 5, column B).
$ws.Range("B5").Value = "Electronics >> Others"

# 2. Restore gridlines/headers display (engine defaults these to hidden
#    on save unless explicitly (re)asserted) so the view stays identical
#    to the source formatting.
$excel.ActiveWindow.DisplayGridlines = $true
$excel.ActiveWindow.DisplayHeadings = $true

# 3. Scroll the window so the top-left visible cell is A1 again, then
#    move the selection/active cell to B11.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B11").Select()
